$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheet1)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 1220
$ws1.Range("F7").Value = 12447
$ws1.Range("F12").Value = 186
$ws1.Range("F13").Value = 12299
$ws1.Range("F14").Value = 4866
$ws1.Range("F15").Value = 4763
$ws1.Range("F19").Value = 98

# Sheet "全部类型" (sheet4)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 1220
$ws4.Range("F9").Value = 12447
$ws4.Range("F14").Value = 186
$ws4.Range("F15").Value = 12299
$ws4.Range("F16").Value = 4866
$ws4.Range("F17").Value = 4763
$ws4.Range("F21").Value = 98
